$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray leftover row (old row 8, which only held "A8=7").
$ws.Rows.Item(8).Delete()

function Set-TextCell($row, $col, $text) {
    # Force the value to be stored as plain text (not auto-parsed into a
    # date/number), then drop back to General formatting so no lingering
    # number format / quote-prefix marker is left on the cell.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
}

# Rewrite the header row.
Set-TextCell 1 1 "#"
Set-TextCell 1 2 "V"
Set-TextCell 1 3 "Place1"
Set-TextCell 1 4 "Place2"
Set-TextCell 1 5 "Q2 (Ton)"
Set-TextCell 1 6 "Date"

# Rewrite the data rows with the updated values.
$ws.Cells.Item(2,1).Value = 1
Set-TextCell 2 2 "Quantum Logistics"
Set-TextCell 2 3 "Shanghai"
Set-TextCell 2 4 "Ningbo"
$ws.Cells.Item(2,5).Value = 1000
Set-TextCell 2 6 "2017/7/1"

$ws.Cells.Item(3,1).Value = 2
Set-TextCell 3 2 "Quantum Logistics"
Set-TextCell 3 3 "Shanghai"
Set-TextCell 3 4 "Nanjing"
$ws.Cells.Item(3,5).Value = 1500
Set-TextCell 3 6 "2017/7/2"

$ws.Cells.Item(4,1).Value = 3
Set-TextCell 4 2 "Quantum Logistics"
Set-TextCell 4 3 "Suzhou"
Set-TextCell 4 4 "Shanghai"
$ws.Cells.Item(4,5).Value = 1000
Set-TextCell 4 6 "2017/7/3"

$ws.Cells.Item(5,1).Value = 4
Set-TextCell 5 2 "Quantum Logistics"
Set-TextCell 5 3 "Shanghai"
Set-TextCell 5 4 "Ningbo"
$ws.Cells.Item(5,5).Value = 1000
Set-TextCell 5 6 "2017/7/4"

$ws.Cells.Item(6,1).ClearContents()
Set-TextCell 6 2 "Quantum Logistics"
Set-TextCell 6 3 "Suzhou"
Set-TextCell 6 4 "Shanghai"
$ws.Cells.Item(6,5).Value = 5000
Set-TextCell 6 6 "2017/7/5"

$ws.Cells.Item(7,1).ClearContents()
Set-TextCell 7 2 "Quantum SCM"
Set-TextCell 7 3 "Guangzhou"
Set-TextCell 7 4 "Shanghai"
$ws.Cells.Item(7,5).Value = 500
Set-TextCell 7 6 "2017/7/6"

# Widen column B (the "V" / vendor column) to fit the longer names.
$ws.Columns.Item(2).ColumnWidth = 16

# Restore the selected cell to match the author's final cursor position.
$ws.Range("F9").Select()
